# Apply "More tweaks to date-of-death states" edit.
$wb = $excel.ActiveWorkbook

# --- Arizona sheet (sheet1) ---
$ws1 = $wb.Worksheets.Item("Arizona")

# Updated daily death counts (column B) for existing rows.
$ws1.Cells.Item(312, 2).Value = 134
$ws1.Cells.Item(314, 2).Value = 114
$ws1.Cells.Item(315, 2).Value = 127
$ws1.Cells.Item(316, 2).Value = 120
$ws1.Cells.Item(317, 2).Value = 119
$ws1.Cells.Item(318, 2).Value = 101
$ws1.Cells.Item(319, 2).Value = 106
$ws1.Cells.Item(320, 2).Value = 112
$ws1.Cells.Item(322, 2).Value = 87
$ws1.Cells.Item(323, 2).Value = 57
$ws1.Cells.Item(324, 2).Value = 68
$ws1.Cells.Item(325, 2).Value = 66
$ws1.Cells.Item(326, 2).Value = 48
$ws1.Cells.Item(327, 2).Value = 30
$ws1.Cells.Item(328, 2).Value = 6

# New row 329: next date, new daily count, cumulative running total formula.
$ws1.Cells.Item(329, 1).Value = 44218
$ws1.Cells.Item(329, 1).NumberFormat = "m/d/yy"
$ws1.Cells.Item(329, 2).Value = 2
$ws1.Cells.Item(329, 3).Formula = "=C328+B329"

# View state updates for the Arizona sheet.
$ws1.Application.ActiveWindow.ScrollRow = 304
$ws1.Range("B318").Select()

# --- New Jersey sheet (sheet2) ---
$ws2 = $wb.Worksheets.Item("New Jersey")
$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollRow = 301
$ws2.Range("C311").Select()

# Re-activate Arizona tab as the active sheet and restore its scroll/selection.
$ws1.Activate()
$ws1.Application.ActiveWindow.ScrollRow = 304
$ws1.Range("B318").Select()

# Workbook window width tweak.
$excel.ActiveWindow.Width = 10000 * 0.75
